$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 260, shifting existing rows 260:390 down to 261:391
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with data
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 45097
$ws.Cells.Item(260, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100112039
$ws.Cells.Item(260, 7).Value = "Ciboulette"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 45
$ws.Cells.Item(260, 11).Value = 7000
$ws.Cells.Item(260, 12).Value = 7000
$ws.Cells.Item(260, 13).Value = 7000
$ws.Cells.Item(260, 14).Value = "`$/docena de atados"
$ws.Cells.Item(260, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(260, 16).Value = 2333
$ws.Cells.Item(260, 17).Value = 3
$ws.Cells.Item(260, 18).Value = "Hortaliza"
